# class 27 Excel File reading
# Adds a new employee record (Moncef / M / Belgas) as row 6 of the
# Employees sheet, right below the existing five rows of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Moncef"
$ws.Range("B6").Value = "M"
$ws.Range("C6").Value = "Belgas"
